# Fix the MATCH() formulas in the "FIT/MAX" analysis sheet.
#
# The sheet used MATCH(value, range) with the match-type argument omitted,
# which defaults to an *approximate* match (match_type = 1) and requires the
# lookup range to be sorted ascending. Since the data ranges are NOT sorted,
# this silently returned wrong results. The fix is to make every MAX-based
# MATCH() call an *exact* match by passing 0 as the third argument - matching
# the MIN-based MATCH() calls elsewhere in the sheet, which already did this
# correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table1 "MAX" column (I2:I8): per-row max across Laplasa..FIT ---
$maxFormula = "=INDEX(Table1[[#Headers],[Laplasa]:[FIT]], MATCH(MAX(Table1[[#This Row],[Laplasa]:[FIT]]), Table1[[#This Row],[Laplasa]:[FIT]], 0))"
for ($row = 2; $row -le 8; $row++) {
    $ws.Range("I$row").Formula = $maxFormula
}

# --- Row 9 ("MAX" summary row, B9:H9): per-column max across rows 2-8 ---
$cols = @("B", "C", "D", "E", "F", "G", "H")
foreach ($col in $cols) {
    $ws.Range($col + "9").Formula = "=INDEX(`$A`$2:`$A`$8, MATCH(MAX(" + $col + "2:" + $col + "8), " + $col + "2:" + $col + "8, 0))"
}

# Reflect the user's last selection/click after fixing the formulas.
[void]$ws.Range("B10").Select()
